# Fruta / hortaliza, semanal
# Insert 2 new daily records at the top of the "Ciruela" block (rows 89-90),
# pushing the existing rows 89-166 down to 91-168.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows above row 89 (existing data shifts down).
$ws.Range("A89:T90").Insert()

# New row 89: Blue Giant
$ws.Range("A89").Value = 10
$ws.Range("B89").Value = "Vega Modelo de Temuco"
$ws.Range("C89").Value = "La Araucanía"
$ws.Range("D89").Value = [DateTime]"2022-01-25"
$ws.Range("D89").NumberFormat = $ws.Range("D91").NumberFormat
$ws.Range("E89").Value = 9
$ws.Range("F89").Value = "Fruta"
$ws.Range("G89").Value = 100103
$ws.Range("H89").Value = "Frutos de hueso (carozo)"
$ws.Range("I89").Value = 100103002
$ws.Range("J89").Value = "Ciruela"
$ws.Range("K89").Value = "Blue Giant"
$ws.Range("L89").Value = "Primera"
$ws.Range("M89").Value = 200
$ws.Range("N89").Value = 18000
$ws.Range("O89").Value = 19000
$ws.Range("P89").Value = 18600
$ws.Range("Q89").Value = "`$/bandeja 18 kilos granel"
$ws.Range("R89").Value = "Región de O'Higgins"
$ws.Range("S89").Value = 1033
$ws.Range("T89").Value = 18

# New row 90: Lemon
$ws.Range("A90").Value = 10
$ws.Range("B90").Value = "Vega Modelo de Temuco"
$ws.Range("C90").Value = "La Araucanía"
$ws.Range("D90").Value = [DateTime]"2022-01-25"
$ws.Range("D90").NumberFormat = $ws.Range("D91").NumberFormat
$ws.Range("E90").Value = 9
$ws.Range("F90").Value = "Fruta"
$ws.Range("G90").Value = 100103
$ws.Range("H90").Value = "Frutos de hueso (carozo)"
$ws.Range("I90").Value = 100103002
$ws.Range("J90").Value = "Ciruela"
$ws.Range("K90").Value = "Lemon"
$ws.Range("L90").Value = "Primera"
$ws.Range("M90").Value = 100
$ws.Range("N90").Value = 15000
$ws.Range("O90").Value = 15000
$ws.Range("P90").Value = 15000
$ws.Range("Q90").Value = "`$/bandeja 18 kilos granel"
$ws.Range("R90").Value = "Región de O'Higgins"
$ws.Range("S90").Value = 833
$ws.Range("T90").Value = 18
